$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows 2-4 (values shifted/changed) and rewrite rows 5
# with the values that used to belong to the old row 4 pattern, then add a
# brand-new row 6 with the values that used to live in the old row 5.

# Row 2: only I2 changes (4 -> 5)
$ws.Range("I2").Value = 5

# Row 3
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0
$ws.Range("H3").Value = 61
$ws.Range("I3").Value = 5

# Row 4
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 9
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 5
$ws.Range("H4").Value = 11
$ws.Range("I4").Value = 5

# Row 5 (new content replacing the previous row 5 values)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 51
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = "train_dim1_2"

# Row 6 (brand new row)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 31
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = "train_dim1_2"
